$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) keeps text formatting so values like "1.002" are not
# auto-converted into numbers/dates by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '29.222.19'
$ws.Range("E2").Value = '  -0.62%  '

# Row 3
$ws.Range("D3").Value = '1.830.25'
$ws.Range("E3").Value = '  -0.71%  '

# Row 5
$ws.Range("D5").Value = '237.25'
$ws.Range("E5").Value = '  -1.24%  '

# Row 6
$ws.Range("D6").Value = '0.6092'
$ws.Range("E6").Value = '  -3.69%  '

# Row 8
$ws.Range("D8").Value = '0.07088'
$ws.Range("E8").Value = '  -5.16%  '

# Row 9
$ws.Range("D9").Value = '0.2815'
$ws.Range("E9").Value = '  -3.14%  '

# Row 10
$ws.Range("D10").Value = '23.85'
$ws.Range("E10").Value = '  -4.96%  '

# Row 11
$ws.Range("D11").Value = '0.07646'
$ws.Range("E11").Value = '  -1.23%  '

# Row 12
$ws.Range("D12").Value = '1.847.31'
$ws.Range("E12").Value = '  +0.12%  '

# Row 13
$ws.Range("D13").Value = '4.805'
$ws.Range("E13").Value = '  -3.73%  '

# Row 14
$ws.Range("D14").Value = '0.6329'
$ws.Range("E14").Value = '  -6.77%  '

# Row 15
$ws.Range("D15").Value = '0.000009995'
$ws.Range("E15").Value = '  -2.48%  '

# Row 16
$ws.Range("D16").Value = '2.067.87'
$ws.Range("E16").Value = '  -1.17%  '

# Row 17
$ws.Range("D17").Value = '79.54'
$ws.Range("E17").Value = '  -3.16%  '

# Row 18
$ws.Range("D18").Value = '5.956'
$ws.Range("E18").Value = '  -5.09%  '

# Row 19
$ws.Range("D19").Value = '29.220.70'
$ws.Range("E19").Value = '  -0.86%  '

# Row 20
$ws.Range("D20").Value = '228.92'

# Row 21
$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D21").Value = '1.002'
$ws.Range("E21").Value = '  +0.26%  '

# Row 22
$ws.Range("B22").Value = 'Avalanche'
$ws.Range("C22").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D22").Value = '11.80'
$ws.Range("E22").Value = '  -4.32%  '

# Row 23
$ws.Range("D23").Value = '7.032'
$ws.Range("E23").Value = '  -5.46%  '

# Row 24
$ws.Range("D24").Value = '1.002'
$ws.Range("E24").Value = '  +0.19%  '

# Row 25
$ws.Range("D25").Value = '155.39'
$ws.Range("E25").Value = '  -1.90%  '

# Row 26
$ws.Range("D26").Value = '8.096'
$ws.Range("E26").Value = '  -4.76%  '

# Row 27
$ws.Range("E27").Value = '  -4.37%  '

# Row 28
$ws.Range("D28").Value = '16.72'
$ws.Range("E28").Value = '  -4.27%  '

# Row 29
$ws.Range("D29").Value = '0.06792'
$ws.Range("E29").Value = '  +2.99%  '

# Row 30
$ws.Range("D30").Value = '1.480'
$ws.Range("E30").Value = '  +3.63%  '

# Row 31
$ws.Range("D31").Value = '1.459'
$ws.Range("E31").Value = '  -1.80%  '

# Row 32
$ws.Range("D32").Value = '3.829'
$ws.Range("E32").Value = '  -6.13%  '

# Row 33
$ws.Range("D33").Value = '3.830'
$ws.Range("E33").Value = '  -5.61%  '

# Row 34
$ws.Range("E34").Value = '  -1.23%  '

# Row 36
$ws.Range("D36").Value = '0.6552'
$ws.Range("E36").Value = '  -6.25%  '

# Row 37
$ws.Range("D37").Value = '2.556'
$ws.Range("E37").Value = '  -0.81%  '

# Row 38
$ws.Range("D38").Value = '1.232.71'
$ws.Range("E38").Value = '  -1.33%  '

# Row 39
$ws.Range("D39").Value = '2.762'
$ws.Range("E39").Value = '  -1.92%  '

# Row 40
$ws.Range("D40").Value = '0.01766'
$ws.Range("E40").Value = '  -5.01%  '

# Row 41
$ws.Range("D41").Value = '6.593'
$ws.Range("E41").Value = '  -2.98%  '

# Row 42
$ws.Range("D42").Value = '0.9232'
$ws.Range("E42").Value = '  -1.20%  '

# Row 43
$ws.Range("D43").Value = '1.001'
$ws.Range("E43").Value = '  +0.20%  '

# Row 44
$ws.Range("D44").Value = '1.984.54'
$ws.Range("E44").Value = '  -0.57%  '

# Row 45
$ws.Range("D45").Value = '100.85'
$ws.Range("E45").Value = '  -0.22%  '

# Row 46
$ws.Range("D46").Value = '63.40'
$ws.Range("E46").Value = '  -3.18%  '

# Row 47
$ws.Range("E47").Value = '  -1.78%  '

# Row 48
$ws.Range("E48").Value = '  -5.30%  '

# Row 49
$ws.Range("D49").Value = '8.525'
$ws.Range("E49").Value = '  -5.89%  '

# Row 50
$ws.Range("B50").Value = 'Algorand'
$ws.Range("C50").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D50").Value = '0.1084'
$ws.Range("E50").Value = '  -5.72%  '

# Row 51
$ws.Range("B51").Value = 'Aptos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D51").Value = '6.524'
$ws.Range("E51").Value = '  -7.71%  '
